$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 5.5
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 1.62
$ws.Range("L2").Value = 1.47
$ws.Range("M2").Value = 2.2
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 25
$ws.Range("P2").Value = 14.5
$ws.Range("Q2").Value = 90
$ws.Range("R2").Value = 50
$ws.Range("S2").Value = 55
$ws.Range("T2").Value = 7.1
$ws.Range("U2").Value = 5.7
$ws.Range("V2").Value = 15.5
$ws.Range("W2").Value = 80
$ws.Range("Y2").Value = 4.45
$ws.Range("Z2").Value = 5.5
$ws.Range("AA2").Value = 7.1
$ws.Range("AB2").Value = 9.75
$ws.Range("AC2").Value = 12.5
$ws.Range("AD2").Value = 28

# Row 8
$ws.Range("N8").Value = 6.3
$ws.Range("P8").Value = 8.75
$ws.Range("S8").Value = 30
$ws.Range("V8").Value = 12
$ws.Range("Y8").Value = 5.7
$ws.Range("Z8").Value = 9.75
$ws.Range("AC8").Value = 19.5

# Row 9
$ws.Range("G9").Value = 1.62
$ws.Range("I9").Value = 5.5
$ws.Range("Q9").Value = 12
$ws.Range("R9").Value = 15
$ws.Range("T9").Value = 8
$ws.Range("Y9").Value = 12
$ws.Range("AA9").Value = 19
$ws.Range("AC9").Value = 51
$ws.Range("AE9").Value = 1.07
$ws.Range("AF9").Value = 9

# Row 15
$ws.Range("H15").Value = 5.1
$ws.Range("I15").Value = 12
$ws.Range("O15").Value = 4.6
$ws.Range("T15").Value = 11.5
$ws.Range("Y15").Value = 23
$ws.Range("Z15").Value = 75
$ws.Range("AA15").Value = 32
$ws.Range("AB15").Value = 350
$ws.Range("AC15").Value = 150

# Row 16
$ws.Range("G16").Value = 2.35
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 2.8
$ws.Range("N16").Value = 6.1
$ws.Range("O16").Value = 9
$ws.Range("P16").Value = 7.9
$ws.Range("Q16").Value = 19
$ws.Range("R16").Value = 16.5
$ws.Range("T16").Value = 7.9
$ws.Range("U16").Value = 5.3
$ws.Range("W16").Value = 60
$ws.Range("Y16").Value = 6.6
$ws.Range("Z16").Value = 11
$ws.Range("AB16").Value = 26
$ws.Range("AC16").Value = 21
$ws.Range("AD16").Value = 30

# Row 17
$ws.Range("J17").Value = 2.1
$ws.Range("K17").Value = 1.7
$ws.Range("AF17").Value = 9

# Row 19
$ws.Range("G19").Value = 1.4
$ws.Range("H19").Value = 3.7
$ws.Range("I19").Value = 7.5
$ws.Range("J19").Value = 2.15
$ws.Range("K19").Value = 1.67
$ws.Range("L19").Value = 1.44
$ws.Range("M19").Value = 2.63
$ws.Range("N19").Value = 5.5
$ws.Range("O19").Value = 5.5
$ws.Range("P19").Value = 9
$ws.Range("Q19").Value = 9
$ws.Range("R19").Value = 15
$ws.Range("S19").Value = 41
$ws.Range("T19").Value = 8
$ws.Range("U19").Value = 8
$ws.Range("V19").Value = 26
$ws.Range("W19").Value = 101
$ws.Range("Y19").Value = 15
$ws.Range("Z19").Value = 41
$ws.Range("AA19").Value = 26
$ws.Range("AB19").Value = 101
$ws.Range("AC19").Value = 67
$ws.Range("AD19").Value = 81
$ws.Range("AE19").Value = 1.08
$ws.Range("AF19").Value = 8
$ws.Range("AG19").Value = 1.36
$ws.Range("AH19").Value = 3
$ws.Range("AI19").Value = 2.5
$ws.Range("AJ19").Value = 1.5

# Row 20
$ws.Range("G20").Value = 1.45
$ws.Range("H20").Value = 3.6
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 2.05
$ws.Range("K20").Value = 1.75
$ws.Range("L20").Value = 1.4
$ws.Range("M20").Value = 2.75
$ws.Range("N20").Value = 5.5
$ws.Range("O20").Value = 6
$ws.Range("P20").Value = 9
$ws.Range("Q20").Value = 9.5
$ws.Range("R20").Value = 13
$ws.Range("S20").Value = 34
$ws.Range("T20").Value = 8.5
$ws.Range("U20").Value = 7.5
$ws.Range("V20").Value = 21
$ws.Range("W20").Value = 81
$ws.Range("Y20").Value = 15
$ws.Range("Z20").Value = 34
$ws.Range("AA20").Value = 23
$ws.Range("AB20").Value = 81
$ws.Range("AC20").Value = 51
$ws.Range("AD20").Value = 67
$ws.Range("AE20").Value = 1.07
$ws.Range("AF20").Value = 8.5
$ws.Range("AG20").Value = 1.33
$ws.Range("AH20").Value = 3.25
$ws.Range("AI20").Value = 2.2
$ws.Range("AJ20").Value = 1.62

# Row 22
$ws.Range("X22").Value = 800

# Row 23
$ws.Range("G23").Value = 2.35
$ws.Range("I23").Value = 3
$ws.Range("O23").Value = 10
$ws.Range("P23").Value = 10
$ws.Range("Q23").Value = 23
$ws.Range("Y23").Value = 7.5
